$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.04910134944039157
$ws.Range("D2").Value = 0.162116166634263
$ws.Range("E2").Value = 0.2537512982573276
$ws.Range("F2").Value = 0.9853616620280334
$ws.Range("G2").Value = 0.002430949633585369
$ws.Range("J2").Value = 0.4283563618163839
$ws.Range("M2").Value = 10.9357286970332
$ws.Range("O2").Value = 2.786437679966582
$ws.Range("C3").Value = 0.04357122867759244
$ws.Range("D3").Value = 0.1560232641946584
$ws.Range("E3").Value = 0.2325967084506857
$ws.Range("F3").Value = 1.029516344211281
$ws.Range("G3").Value = 0.002437072292602609
$ws.Range("J3").Value = 0.3848457021135516
$ws.Range("M3").Value = 9.612703468246707
$ws.Range("O3").Value = 2.862946529673025
$ws.Range("C4").Value = 0.04019728211902418
$ws.Range("D4").Value = 0.1523631980409021
$ws.Range("E4").Value = 0.2197763238064852
$ws.Range("F4").Value = 1.058762334134251
$ws.Range("G4").Value = 0.002440993109030101
$ws.Range("J4").Value = 0.3584113181624957
$ws.Range("M4").Value = 8.797400819332893
$ws.Range("O4").Value = 2.915534710963414
$ws.Range("C5").Value = 0.03882757843490481
$ws.Range("D5").Value = 0.1508918880495997
$ws.Range("E5").Value = 0.2145927303900308
$ws.Range("F5").Value = 1.071207773722179
$ws.Range("G5").Value = 0.002442631707241562
$ws.Range("J5").Value = 0.3477061187714412
$ws.Range("M5").Value = 8.464380202270775
$ws.Range("O5").Value = 2.938350607595339
$ws.Range("C6").Value = 0.03860044845109201
$ws.Range("D6").Value = 0.1506487936340619
$ws.Range("E6").Value = 0.2137344202183655
$ws.Range("F6").Value = 1.073305915187952
$ws.Range("G6").Value = 0.002442906268184312
$ws.Range("J6").Value = 0.3459324734913025
$ws.Range("M6").Value = 8.409034484731251
$ws.Range("O6").Value = 2.942222134642918
$ws.Range("C7").Value = 0.04017878901555605
$ws.Range("D7").Value = 0.1523432738649859
$ws.Range("E7").Value = 0.2197062527792752
$ws.Range("F7").Value = 1.058928054289758
$ws.Range("G7").Value = 0.002441015041877974
$ws.Range("J7").Value = 0.3582666774217387
$ws.Range("M7").Value = 8.7929127789983
$ws.Range("O7").Value = 2.915836836265953
$ws.Range("C8").Value = 0.04718995276459736
$ws.Range("D8").Value = 0.159998405910315
$ws.Range("E8").Value = 0.2464211977394939
$ws.Range("F8").Value = 1.000137451127415
$ws.Range("G8").Value = 0.002433027352024647
$ws.Range("J8").Value = 0.4132931576789929
$ws.Range("M8").Value = 10.48014164963467
$ws.Range("O8").Value = 2.811638000453399
$ws.Range("C9").Value = 0.0611207289749558
$ws.Range("D9").Value = 0.1756633327369883
$ws.Range("E9").Value = 0.3002280680070442
$ws.Range("F9").Value = 0.9022112971321121
$ws.Range("G9").Value = 0.002418633901242167
$ws.Range("J9").Value = 0.5236194466843642
$ws.Range("M9").Value = 13.76723892034278
$ws.Range("O9").Value = 2.652973941649918
$ws.Range("C10").Value = 0.07148288808997449
$ws.Range("D10").Value = 0.1875893170603717
$ws.Range("E10").Value = 0.3407515770546041
$ws.Range("F10").Value = 0.8414379404937051
$ws.Range("G10").Value = 0.002408817628314487
$ws.Range("J10").Value = 0.6064397348642956
$ws.Range("M10").Value = 16.17239597373094
$ws.Range("O10").Value = 2.565843596840836
$ws.Range("C11").Value = 0.07622826814424855
$ws.Range("D11").Value = 0.1931102807944853
$ws.Range("E11").Value = 0.359433002896921
$ws.Range("F11").Value = 0.8163471564696962
$ws.Range("G11").Value = 0.002404513089326896
$ws.Range("J11").Value = 0.6445695251636323
$ws.Range("M11").Value = 17.2652529995803
$ws.Range("O11").Value = 2.532952549403632
$ws.Range("C12").Value = 0.07803005823552667
$ws.Range("D12").Value = 0.1952151399285071
$ws.Range("E12").Value = 0.3665454014896312
$ws.Range("F12").Value = 0.8072245450992881
$ws.Range("G12").Value = 0.002402905922416597
$ws.Range("J12").Value = 0.6590797070934116
$ws.Range("M12").Value = 17.67897970287771
$ws.Range("O12").Value = 2.521497563261818
$ws.Range("C13").Value = 0.07764179172473007
$ws.Range("D13").Value = 0.1947611825851823
$ws.Range("E13").Value = 0.3650118787719521
$ws.Range("F13").Value = 0.8091722362028833
$ws.Range("G13").Value = 0.00240325104204001
$ws.Range("J13").Value = 0.655951412055316
$ws.Range("M13").Value = 17.58988021726321
$ws.Range("O13").Value = 2.523919628716612
$ws.Range("C14").Value = 0.07637640432666615
$ws.Range("D14").Value = 0.193283161873552
$ws.Range("E14").Value = 0.3600173639745066
$ws.Range("F14").Value = 0.8155889779566223
$ws.Range("G14").Value = 0.002404380410090616
$ws.Range("J14").Value = 0.6457618231236495
$ws.Range("M14").Value = 17.29929243447788
$ws.Range("O14").Value = 2.531989920760907
$ws.Range("C15").Value = 0.07560195378496815
$ws.Range("D15").Value = 0.1923796926779175
$ws.Range("E15").Value = 0.3569631284868962
$ws.Range("F15").Value = 0.8195690781147107
$ws.Range("G15").Value = 0.002405075150578692
$ws.Range("J15").Value = 0.639529868025221
$ws.Range("M15").Value = 17.12128620227514
$ws.Range("O15").Value = 2.537064365069227
$ws.Range("C16").Value = 0.07117342522747094
$ws.Range("D16").Value = 0.1872304709428931
$ws.Range("E16").Value = 0.3395358932150572
$ws.Range("F16").Value = 0.8431300587972643
$ws.Range("G16").Value = 0.002409102155313234
$ws.Range("J16").Value = 0.6039575036884628
$ws.Range("M16").Value = 16.10095538832451
$ws.Range("O16").Value = 2.56813139670956
$ws.Range("C17").Value = 0.06846494855633978
$ws.Range("D17").Value = 0.1840964024605682
$ws.Range("E17").Value = 0.3289099337243471
$ws.Range("F17").Value = 0.8582463659384558
$ws.Range("G17").Value = 0.002411613619424163
$ws.Range("J17").Value = 0.5822554779055906
$ws.Range("M17").Value = 15.47473804671301
$ws.Range("O17").Value = 2.588938834743772
$ws.Range("C18").Value = 0.06691005990548149
$ws.Range("D18").Value = 0.182302755934046
$ws.Range("E18").Value = 0.3228212506648589
$ws.Range("F18").Value = 0.8671806939140936
$ws.Range("G18").Value = 0.002413073311698665
$ws.Range("J18").Value = 0.5698154410074778
$ws.Range("M18").Value = 15.1144313218611
$ws.Range("O18").Value = 2.601539741913911
$ws.Range("C19").Value = 0.06638410068660505
$ws.Range("D19").Value = 0.1816969895673424
$ws.Range("E19").Value = 0.3207636159141742
$ws.Range("F19").Value = 0.8702465557505548
$ws.Range("G19").Value = 0.002413570151264189
$ws.Range("J19").Value = 0.565610566027658
$ws.Range("M19").Value = 14.99241458575534
$ws.Range("O19").Value = 2.605914005488643
$ws.Range("C20").Value = 0.06875296267138253
$ws.Range("D20").Value = 0.1844290964286586
$ws.Range("E20").Value = 0.3300386774008359
$ws.Range("F20").Value = 0.8566123072353733
$ws.Range("G20").Value = 0.002411344702323322
$ws.Range("J20").Value = 0.5845612664170972
$ws.Range("M20").Value = 15.54141222556547
$ws.Range("O20").Value = 2.58665811585638
$ws.Range("C21").Value = 0.07674794627662607
$ws.Range("D21").Value = 0.193716903701926
$ws.Range("E21").Value = 0.3614833168744838
$ws.Range("F21").Value = 0.813693855594039
$ws.Range("G21").Value = 0.002404048068892886
$ws.Range("J21").Value = 0.648752769320339
$ws.Range("M21").Value = 17.3846476549358
$ws.Range("O21").Value = 2.52959209223485
$ws.Range("C22").Value = 0.0820013469242582
$ws.Range("D22").Value = 0.1998699610826975
$ws.Range("E22").Value = 0.3822579181084649
$ws.Range("F22").Value = 0.7878571050326997
$ws.Range("G22").Value = 0.002399412458807486
$ws.Range("J22").Value = 0.6911240117993032
$ws.Range("M22").Value = 18.58867652350335
$ws.Range("O22").Value = 2.498140302694679
$ws.Range("C23").Value = 0.07919483229939317
$ws.Range("D23").Value = 0.1965782191145706
$ws.Range("E23").Value = 0.371148736645182
$ws.Range("F23").Value = 0.8014403791167766
$ws.Range("G23").Value = 0.00240187447970057
$ws.Range("J23").Value = 0.6684693392201382
$ws.Range("M23").Value = 17.94609775963994
$ws.Range("O23").Value = 2.514381789497321
$ws.Range("C24").Value = 0.06862274450045902
$ws.Range("D24").Value = 0.1842786600651891
$ws.Range("E24").Value = 0.3295283092800503
$ws.Range("F24").Value = 0.8573503067617594
$ws.Range("G24").Value = 0.002411466230442884
$ws.Range("J24").Value = 0.5835187054151731
$ws.Range("M24").Value = 15.51126972520802
$ws.Range("O24").Value = 2.587687242564755
$ws.Range("C25").Value = 0.05733095622518647
$ws.Range("D25").Value = 0.171353842134593
$ws.Range("E25").Value = 0.285507522792031
$ws.Range("F25").Value = 0.9267841457537145
$ws.Range("G25").Value = 0.002422393317935045
$ws.Range("J25").Value = 0.4934857321182164
$ws.Range("M25").Value = 12.87998317466628
$ws.Range("O25").Value = 2.690844331968322
